# Auto-generated edit script: apply cached-value updates to Anima Profits sheets
# (workbook has no formulas; all target cells hold static numeric literals)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 82
$ws.Range("I11").Value = 82
$ws.Range("K11").Value = 82
$ws.Range("M11").Value = 58
$ws.Range("H15").Value = 634.24
$ws.Range("I15").Value = 634.24
$ws.Range("K15").Value = 1902.72
$ws.Range("M15").Value = -1733.72
$ws.Range("H75").Value = 36016.668
$ws.Range("J75").Value = 36016.668
$ws.Range("L75").Value = 36016.668
$ws.Range("N75").Value = -37888.668
$ws.Range("H78").Value = 36016.668
$ws.Range("J78").Value = 36016.668
$ws.Range("L78").Value = 108050.004
$ws.Range("N78").Value = -117410.004
$ws.Range("H132").Value = 2234.4558
$ws.Range("I132").Value = 2135.9314
$ws.Range("J132").Value = 3433.1667
$ws.Range("K132").Value = 6407.7942
$ws.Range("L132").Value = 10299.5001
$ws.Range("M132").Value = -3877.7942
$ws.Range("N132").Value = -15359.5001
$ws.Range("H137").Value = 1249.5231
$ws.Range("I137").Value = 1063.238
$ws.Range("J137").Value = 1589.6957
$ws.Range("K137").Value = 3189.714
$ws.Range("L137").Value = 4769.0871
$ws.Range("M137").Value = -639.7139999999999
$ws.Range("N137").Value = -9869.087100000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 180.2
$ws.Range("I5").Value = 150.5
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 150.5
$ws.Range("L5").Value = 200
$ws.Range("M5").Value = -38.5
$ws.Range("N5").Value = -424
$ws.Range("H11").Value = 2000476
$ws.Range("I11").Value = 2666966.8
$ws.Range("J11").Value = 1004
$ws.Range("K11").Value = 2666966.8
$ws.Range("L11").Value = 1004
$ws.Range("M11").Value = -2666822.8
$ws.Range("N11").Value = -1292
$ws.Range("H61").Value = 7938651
$ws.Range("I61").Value = 9805525
$ws.Range("J61").Value = 4437.375
$ws.Range("K61").Value = 9805525
$ws.Range("L61").Value = 4437.375
$ws.Range("M61").Value = -9805313
$ws.Range("N61").Value = -4861.375
$ws.Range("H74").Value = 1178.35
$ws.Range("I74").Value = 980.7222
$ws.Range("J74").Value = 2957
$ws.Range("K74").Value = 980.7222
$ws.Range("L74").Value = 2957
$ws.Range("M74").Value = -106.7222
$ws.Range("N74").Value = -4705
$ws.Range("H77").Value = 1178.35
$ws.Range("I77").Value = 980.7222
$ws.Range("J77").Value = 2957
$ws.Range("K77").Value = 4903.611
$ws.Range("L77").Value = 14785
$ws.Range("M77").Value = -535.6109999999999
$ws.Range("N77").Value = -23521
$ws.Range("H97").Value = 1077.9524
$ws.Range("I97").Value = 1116.5385
$ws.Range("J97").Value = 1015.25
$ws.Range("K97").Value = 1116.5385
$ws.Range("L97").Value = 1015.25
$ws.Range("M97").Value = -620.5385000000001
$ws.Range("N97").Value = -2007.25
$ws.Range("H102").Value = 2332
$ws.Range("I102").Value = 2045.7142
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 2045.7142
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -423.7141999999999
$ws.Range("N102").Value = -6244
$ws.Range("H136").Value = 7938651
$ws.Range("I136").Value = 9805525
$ws.Range("J136").Value = 4437.375
$ws.Range("K136").Value = 29416575
$ws.Range("L136").Value = 13312.125
$ws.Range("M136").Value = -29414025
$ws.Range("N136").Value = -18412.125

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 180.2
$ws.Range("I4").Value = 150.5
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 150.5
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = -35.5
$ws.Range("N4").Value = -430
$ws.Range("H11").Value = 4900
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 4900
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 4900
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -5180
$ws.Range("H12").Value = 419.75
$ws.Range("I12").Value = 419.75
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 419.75
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -251.75
$ws.Range("N12").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1678.8918
$ws.Range("I122").Value = 1362.875
$ws.Range("J122").Value = 1919.6666
$ws.Range("K122").Value = 4088.625
$ws.Range("L122").Value = 5758.9998
$ws.Range("M122").Value = -1638.625
$ws.Range("N122").Value = -10658.9998
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800
$ws.Range("H132").Value = 4387308
$ws.Range("I132").Value = 1249.5161
$ws.Range("J132").Value = 23811282
$ws.Range("K132").Value = 3748.5483
$ws.Range("L132").Value = 71433846
$ws.Range("M132").Value = -1218.5483
$ws.Range("N132").Value = -71438906

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 258.52942
$ws.Range("I6").Value = 50.4
$ws.Range("J6").Value = 555.8570999999999
$ws.Range("K6").Value = 151.2
$ws.Range("L6").Value = 1667.5713
$ws.Range("M6").Value = -38.19999999999999
$ws.Range("N6").Value = -1893.5713
$ws.Range("H13").Value = 845.8570999999999
$ws.Range("I13").Value = 30
$ws.Range("J13").Value = 981.8333
$ws.Range("K13").Value = 90
$ws.Range("L13").Value = 2945.4999
$ws.Range("M13").Value = 78
$ws.Range("N13").Value = -3281.4999
$ws.Range("H29").Value = 30303504
$ws.Range("J29").Value = 37037572
$ws.Range("L29").Value = 111112716
$ws.Range("N29").Value = -111113270
$ws.Range("H75").Value = 515
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 515
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H87").Value = 7000
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 7000
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H93").Value = 3644.2666
$ws.Range("I93").Value = 819
$ws.Range("K93").Value = 2457
$ws.Range("M93").Value = -585
$ws.Range("H98").Value = 473
$ws.Range("I98").Value = 326.5
$ws.Range("J98").Value = 598.5714
$ws.Range("K98").Value = 979.5
$ws.Range("L98").Value = 1795.7142
$ws.Range("M98").Value = 518.5
$ws.Range("N98").Value = -4791.7142
$ws.Range("H129").Value = 1064.0286
$ws.Range("I129").Value = 474.91666
$ws.Range("J129").Value = 1371.3914
$ws.Range("K129").Value = 1424.74998
$ws.Range("L129").Value = 4114.174199999999
$ws.Range("M129").Value = 3575.25002
$ws.Range("N129").Value = -14114.1742
$ws.Range("H131").Value = 3076.4717
$ws.Range("J131").Value = 4482.971
$ws.Range("L131").Value = 13448.913
$ws.Range("N131").Value = -23528.913
$ws.Range("H137").Value = 4391635.5
$ws.Range("I137").Value = 9811712
$ws.Range("J137").Value = 3954.7144
$ws.Range("K137").Value = 29435136
$ws.Range("L137").Value = 11864.1432
$ws.Range("M137").Value = -29430036
$ws.Range("N137").Value = -22064.1432
$ws.Range("H139").Value = 3004.2778
$ws.Range("I139").Value = 1958.4615
$ws.Range("J139").Value = 3595.3914
$ws.Range("K139").Value = 5875.3845
$ws.Range("L139").Value = 10786.1742
$ws.Range("M139").Value = -735.3845000000001
$ws.Range("N139").Value = -21066.1742

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 38000
$ws.Range("J32").Value = 38000
$ws.Range("L32").Value = 38000
$ws.Range("N32").Value = -38592

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 766.6667
$ws.Range("I46").Value = 566.6667
$ws.Range("J46").Value = 1166.6666
$ws.Range("K46").Value = 566.6667
$ws.Range("L46").Value = 1166.6666
$ws.Range("M46").Value = -378.6667
$ws.Range("N46").Value = -1542.6666
$ws.Range("H55").Value = 353.02856
$ws.Range("I55").Value = 238.38889
$ws.Range("J55").Value = 474.41177
$ws.Range("K55").Value = 238.38889
$ws.Range("L55").Value = 474.41177
$ws.Range("M55").Value = -65.38889
$ws.Range("N55").Value = -820.4117699999999
$ws.Range("H136").Value = 3969734.2
$ws.Range("I136").Value = 1299.4814
$ws.Range("J136").Value = 11112917
$ws.Range("K136").Value = 3898.4442
$ws.Range("L136").Value = 33338751
$ws.Range("M136").Value = -1348.4442
$ws.Range("N136").Value = -33343851

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2196.02
$ws.Range("I122").Value = 1908.1177
$ws.Range("J122").Value = 2807.8125
$ws.Range("K122").Value = 5724.3531
$ws.Range("L122").Value = 8423.4375
$ws.Range("M122").Value = -3274.3531
$ws.Range("N122").Value = -13323.4375
$ws.Range("H132").Value = 5030035.5
$ws.Range("I132").Value = 1436.7028
$ws.Range("J132").Value = 13889948
$ws.Range("K132").Value = 4310.1084
$ws.Range("L132").Value = 41669844
$ws.Range("M132").Value = -1780.1084
$ws.Range("N132").Value = -41674904
$ws.Range("H133").Value = 41000
$ws.Range("J133").Value = 41000
$ws.Range("L133").Value = 41000
$ws.Range("N133").Value = -51120
$ws.Range("H136").Value = 680.5599999999999
$ws.Range("I136").Value = 592.1686999999999
$ws.Range("J136").Value = 1112.1177
$ws.Range("K136").Value = 1776.5061
$ws.Range("L136").Value = 3336.3531
$ws.Range("M136").Value = 773.4939000000002
$ws.Range("N136").Value = -8436.3531

Write-Host "Applied 245 cell edits across Sheets ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR"
